$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold plain decimal numbers (no thousands separators) that
# Excel would otherwise auto-convert to a Number type on assignment. Force
# them to stay Text (matching the original inlineStr cell type) by setting
# the number format to Text ("@") before writing the new value.
$textPriceCells = @("D5","D8","D16","D18","D19","D23","D25","D28","D38","D40","D43","D44","D48","D51")
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.864.81"
$ws.Range("E2").Value = "  -0.27%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.633.33"
$ws.Range("E3").Value = "  -0.22%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
$ws.Range("D5").Value = "211.57"
$ws.Range("E5").Value = "  -0.42%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.75%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.09%  "

# Row 8 - Solana
$ws.Range("D8").Value = "23.37"
$ws.Range("E8").Value = "  +0.16%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.04%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.15%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.864.60"
$ws.Range("E12").Value = "  -0.24%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.644.76"
$ws.Range("E13").Value = "  +0.57%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -1.21%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -1.28%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "65.30"
$ws.Range("E16").Value = "  +0.04%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.886.22"
$ws.Range("E17").Value = "  -0.22%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "228.75"
$ws.Range("E18").Value = "  -1.00%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "7.68"
$ws.Range("E19").Value = "  +1.98%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -0.21%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.03%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.92%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "10.03"
$ws.Range("E23").Value = "  -3.78%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +0.18%  "

# Row 25 - Monero
$ws.Range("D25").Value = "155.15"
$ws.Range("E25").Value = "  +0.81%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  -1.53%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +0.07%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "15.54"
$ws.Range("E28").Value = "  -0.58%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.07%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.36%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.11%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.90%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +1.16%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.393.68"
$ws.Range("E34").Value = "  -1.06%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +0.92%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  +6.73%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  -0.64%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.0171"
$ws.Range("E38").Value = "  +0.60%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -0.66%  "

# Row 40 - ARBITRUM
$ws.Range("D40").Value = "0.849"
$ws.Range("E40").Value = "  -3.03%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.08%  "

# Row 42 - WEMIXToken
$ws.Range("E42").Value = "  -0.94%  "

# Row 43 - Aave
$ws.Range("D43").Value = "65.85"
$ws.Range("E43").Value = "  -1.96%  "

# Row 44 - RenderToken
$ws.Range("D44").Value = "1.83"
$ws.Range("E44").Value = "  +0.32%  "

# Row 45 - FraxShare
$ws.Range("E45").Value = "  -1.52%  "

# Row 46 - RocketPoolETH (price only; volume unchanged)
$ws.Range("D46").Value = "1.773.40"

# Row 47 - MXToken
$ws.Range("E47").Value = "  -2.80%  "

# Row 48 - Quant
$ws.Range("D48").Value = "88.69"
$ws.Range("E48").Value = "  +0.81%  "

# Row 49 - Algorand
$ws.Range("E49").Value = "  +1.96%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  -0.17%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "7.61"
$ws.Range("E51").Value = "  +0.36%  "
